{"js": "// Insert a new \"Author\"-styled paragraph right after the existing\n// \"Edison Achalma\" author paragraph, containing the affiliation line:\n// \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal\n// de Huamanga\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the paragraph styled \"Author\" whose text is exactly\n// \"Edison Achalma\" (the byline paragraph just below the \"Editar: Editar\"\n// heading \u2014 not the later \"Nota de Autores\" occurrences of the name).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Author\" && p.text.trim() === \"Edison Achalma\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"Edison Achalma\" Author paragraph.');\n}\n\n// NOTE: calling `target.insertParagraph(text, \"After\")` directly on this\n// Paragraph object destroys the paragraph's existing text in this runtime.\n// Inserting a paragraph break + text through the end-of-range via\n// insertText keeps the original paragraph intact and lets the new\n// paragraph correctly inherit the \"Author\" style from the paragraph mark.\nconst endRange = target.getRange(\"End\");\nendRange.insertText(\n  \"\\rEscuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "# Insert a new \"Author\"-styled paragraph right after the existing\n# \"Edison Achalma\" author paragraph, containing the affiliation line:\n# \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal\n# de Huamanga\".\n$d = $word.ActiveDocument\n\n# Locate the paragraph styled \"Author\" whose text is exactly\n# \"Edison Achalma\" (the byline paragraph just below the \"Editar: Editar\"\n# heading \u2014 not the later \"Nota de Autores\" occurrences of the name).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Author\" -and $p.Range.Text.Trim() -eq \"Edison Achalma\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Edison Achalma' Author paragraph.\"\n}\n\n# NOTE: $target.Range.InsertParagraphAfter() (and InsertParagraphBefore /\n# InsertParagraph) drop the paragraph's existing text in this runtime.\n# Appending \"<CR>\" + the new text via InsertAfter keeps the original\n# paragraph intact and the new paragraph correctly inherits the \"Author\"\n# style from the paragraph mark.\n$target.Range.InsertAfter([char]13 + \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\")\n"}
